# Recurso modificado por presencia de ecuaciones / Archivos actualizados
#
# - Renumera la columna "Núm." de la tabla de imágenes (F1..F7 -> IMG01..IMG07)
# - Corrige la referencia rota (#REF!) en las fórmulas de las columnas F y H
#   de la hoja "Solicitud gráfica" para que usen $C$7 (el código de guión/recurso)
# - Mueve la selección activa de A16 a A12

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Renumerar columna A (Núm.) filas 10-16 ---
$ws.Range("A10").Value = "IMG01"
$ws.Range("A11").Value = "IMG02"
$ws.Range("A12").Value = "IMG03"
$ws.Range("A13").Value = "IMG04"
$ws.Range("A14").Value = "IMG05"
$ws.Range("A15").Value = "IMG06"
$ws.Range("A16").Value = "IMG07"

# --- 2. Corregir fórmulas columna F (nombre imagen "small"/normal) filas 10-16 ---
$ws.Range("F10").Formula = '=IF(OR(B10<>"",J10<>""),CONCATENATE($C$7,"_",$A10,IF($G$4="Cuaderno de Estudio","_small",CONCATENATE(IF(I10="","","n"),IF(LEFT($G$5,1)="F",".jpg",".png")))),"")'
$ws.Range("F11:F16").Formula = '=IF(OR(B11<>"",J11<>""),CONCATENATE($C$7,"_",$A11,IF($G$4="Cuaderno de Estudio","_small",CONCATENATE(IF(I11="","","n"),IF(LEFT($G$5,1)="F",".jpg",".png")))),"")'

# --- 3. Corregir fórmulas columna H (nombre imagen "zoom"/ampliada) filas 10-108 ---
$ws.Range("H10").Formula = '=IF(AND(I10<>"",I10<>0),IF(OR(B10<>"",J10<>""),CONCATENATE($C$7,"_",$A10,IF($G$4="Cuaderno de Estudio","_zoom",CONCATENATE("a",IF(LEFT($G$5,1)="F",".jpg",".png")))),""),"")'
$ws.Range("H11:H108").Formula = '=IF(AND(I11<>"",I11<>0),IF(OR(B11<>"",J11<>""),CONCATENATE($C$7,"_",$A11,IF($G$4="Cuaderno de Estudio","_zoom",CONCATENATE("a",IF(LEFT($G$5,1)="F",".jpg",".png")))),""),"")'

# --- 4. Mover la celda activa seleccionada de A16 a A12 ---
$ws.Range("A12").Select()
